# Binance ETHUSDT daily-candle data sheet.
# This commit (2020-11-15) refreshes the last existing candle row (165,
# 2020-11-11) with a later-fetched close snapshot and appends four new
# daily candles (166-169, 2020-11-12 .. 2020-11-15), each with its
# recomputed moving-average / EMA / MACD (DIF/DEM/OSC) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    # Preserve the source file's convention of storing these numeric-looking
    # values (fixed decimal strings) as literal text, not as a Number type.
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
}

function Set-NumCell($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

# --- 1) Update existing row 165: candle was updated with a later close snapshot,
#        shifting derived OHLCV + moving-average / EMA / MACD columns. ---
Set-TextCell 165 3 "476.25000000"
Set-TextCell 165 5 "463.09000000"
Set-TextCell 165 6 "1205782.84262000"
Set-TextCell 165 8 "559228791.89035090"
Set-NumCell 165 9 560360
Set-TextCell 165 10 "626267.49881000"
Set-TextCell 165 11 "290487786.95615110"
Set-NumCell 165 13 463.0899999999996
Set-NumCell 165 14 456.7149999999997
Set-NumCell 165 15 445.6942857142858
Set-NumCell 165 16 415.6153333333334
Set-NumCell 165 17 401.9596666666665
Set-NumCell 165 18 463.09
Set-NumCell 165 19 458.3276506874824
Set-NumCell 165 20 430.9713775808049
Set-NumCell 165 21 411.4122145110911
Set-NumCell 165 22 19.55916306971375
Set-NumCell 165 23 14.0535880226521
Set-NumCell 165 24 5.505575047061649

# --- 2) Append 4 new daily rows (166-169) for 2020-11-12 .. 2020-11-15 ---
# Row 166
$ws.Cells.Item(165, 1).Copy($ws.Cells.Item(166, 1))
Set-NumCell 166 1 164
Set-TextCell 166 2 "463.09000000"
Set-TextCell 166 3 "470.00000000"
Set-TextCell 166 4 "451.20000000"
Set-TextCell 166 5 "462.39000000"
Set-TextCell 166 6 "990776.62912000"
Set-NumCell 166 7 1605225599999
Set-TextCell 166 8 "455971309.21685460"
Set-NumCell 166 9 411840
Set-TextCell 166 10 "471577.98126000"
Set-TextCell 166 11 "217115755.31874190"
Set-TextCell 166 12 "2020-11-12 08:00:00"
Set-NumCell 166 13 462.3899999999996
Set-NumCell 166 14 462.7399999999997
Set-NumCell 166 15 452.2228571428572
Set-NumCell 166 16 420.5593333333334
Set-NumCell 166 17 404.6713333333331
Set-NumCell 166 18 462.39
Set-NumCell 166 19 461.0358835624941
Set-NumCell 166 20 435.8050117991477
Set-NumCell 166 21 415.1883583074846
Set-NumCell 166 22 20.61665349166316
Set-NumCell 166 23 15.36620111645431
Set-NumCell 166 24 5.250452375208846

# Row 167
$ws.Cells.Item(166, 1).Copy($ws.Cells.Item(167, 1))
Set-NumCell 167 1 165
Set-TextCell 167 2 "462.48000000"
Set-TextCell 167 3 "478.01000000"
Set-TextCell 167 4 "457.12000000"
Set-TextCell 167 5 "476.43000000"
Set-TextCell 167 6 "976665.09752000"
Set-NumCell 167 7 1605311999999
Set-TextCell 167 8 "456232937.20613030"
Set-NumCell 167 9 427363
Set-TextCell 167 10 "509344.22798000"
Set-TextCell 167 11 "238031863.05169850"
Set-TextCell 167 12 "2020-11-13 08:00:00"
Set-NumCell 167 13 476.4299999999996
Set-NumCell 167 14 469.4099999999997
Set-NumCell 167 15 455.1542857142858
Set-NumCell 167 16 426.5126666666667
Set-NumCell 167 17 407.9293333333331
Set-NumCell 167 18 476.43
Set-NumCell 167 19 471.2986278541647
Set-NumCell 167 20 442.0550099838998
Set-NumCell 167 21 419.7247890448951
Set-NumCell 167 22 22.33022093900473
Set-NumCell 167 23 16.7590050809644
Set-NumCell 167 24 5.571215858040336

# Row 168
$ws.Cells.Item(167, 1).Copy($ws.Cells.Item(168, 1))
Set-NumCell 168 1 166
Set-TextCell 168 2 "476.42000000"
Set-TextCell 168 3 "477.47000000"
Set-TextCell 168 4 "452.00000000"
Set-TextCell 168 5 "460.89000000"
Set-TextCell 168 6 "735252.78540000"
Set-NumCell 168 7 1605398399999
Set-TextCell 168 8 "340142257.19369760"
Set-NumCell 168 9 350772
Set-TextCell 168 10 "347000.03929000"
Set-TextCell 168 11 "160634962.00481430"
Set-TextCell 168 12 "2020-11-14 08:00:00"
Set-NumCell 168 13 460.8899999999996
Set-NumCell 168 14 468.6599999999998
Set-NumCell 168 15 458.8228571428572
Set-NumCell 168 16 431.7393333333335
Set-NumCell 168 17 410.7043333333331
Set-NumCell 168 18 460.89
Set-NumCell 168 19 464.3595426180549
Set-NumCell 168 20 444.9527007556097
Set-NumCell 168 21 422.7740719183204
Set-NumCell 168 22 22.17862883728935
Set-NumCell 168 23 17.84292983222939
Set-NumCell 168 24 4.335699005059965

# Row 169
$ws.Cells.Item(168, 1).Copy($ws.Cells.Item(169, 1))
Set-NumCell 169 1 167
Set-TextCell 169 2 "460.90000000"
Set-TextCell 169 3 "462.89000000"
Set-TextCell 169 4 "456.51000000"
Set-TextCell 169 5 "461.24000000"
Set-TextCell 169 6 "128673.84301000"
Set-NumCell 169 7 1605484799999
Set-TextCell 169 8 "59155584.09075830"
Set-NumCell 169 9 66228
Set-TextCell 169 10 "67273.91394000"
Set-TextCell 169 11 "30926653.60773300"
Set-TextCell 169 12 "2020-11-15 08:00:00"
Set-NumCell 169 13 461.2399999999997
Set-NumCell 169 14 461.0649999999997
Set-NumCell 169 15 459.8142857142857
Set-NumCell 169 16 436.7246666666668
Set-NumCell 169 17 413.8956666666664
Set-NumCell 169 18 461.24
Set-NumCell 169 19 462.2798475393517
Set-NumCell 169 20 447.4584391009021
Set-NumCell 169 21 425.6234068355566
Set-NumCell 169 22 21.83503226534543
Set-NumCell 169 23 18.6413503188526
Set-NumCell 169 24 3.193681946492834

